# Apply "New report generation method" changes to SalesReport.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the report title text
$ws.Range("A1").Value = "Sales for the last year (365 days)"

# Update the data row (row 3) with the new computed totals
$ws.Range("A3").Value = 49.75
$ws.Range("B3").Value = 60.97
$ws.Range("C3").Value = 16.5
$ws.Range("G3").Value = 74.75
